# Daily refresh of the tracker sheet.
# For every data row:
#   end date   = F (start date, yyyyMMdd) + D (total days)
#   remaining  = end date - "today" (2026-02-06, the day this refresh represents)
# If remaining has run out (<= 0) the cycle restarts: F becomes "today" and
# E is reset back to the full D day count. Otherwise E is simply updated to
# the freshly computed remaining-day count and F is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = [datetime]::ParseExact("20260206", "yyyyMMdd", $null)
$todaySerial = $today.ToOADate()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {

    $days = $ws.Cells.Item($row, 4).Value2
    $startRaw = $ws.Cells.Item($row, 6).Value2

    if ($days -eq $null -or $startRaw -eq $null) { continue }

    $parsed = $true
    try {
        $startDate = [datetime]::ParseExact([string]$startRaw, "yyyyMMdd", $null)
    } catch {
        $parsed = $false
    }
    if (-not $parsed) { continue }

    $endDate = $startDate.AddDays($days)
    $remaining = [int]($endDate.ToOADate() - $todaySerial)

    if ($remaining -le 0) {
        $ws.Cells.Item($row, 5).Value = $days
        $ws.Cells.Item($row, 6).Value = [int]$today.ToString("yyyyMMdd")
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining
    }
}
